## Diagrama de Gantt - add "Logótipo" milestone row (new row 16) to the
## "Marcos435" table, shifting the existing rows 16-25 down to 17-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$lo = $ws.ListObjects.Item("Marcos435")

## 1. Insert a blank worksheet row at row 16 - this pushes the current
##    row 16 ("2a Entrega") and everything below it down by one row.
$ws.Rows.Item(16).Insert()

## 2. Clone the formatting (styles + row height) of row 15 ("Apresentação")
##    into the freshly inserted row 16, since it is a normal (non-header)
##    table data row with the same look & feel as the new "Logótipo" row.
$ws.Range("A15:BM15").Copy($ws.Range("A16"))
$ws.Rows.Item(16).RowHeight = 40.15

## 3. Populate the new row 16 with the "Logótipo" milestone data.
$ws.Range("B16").Value = "      Logótipo"
$ws.Range("C16").Value = "Risco Baixo"
$ws.Range("D16").Value = "Fausto"
$ws.Range("E16").Value2 = 1
$ws.Range("F16").Value2 = 45584
$ws.Range("G16").Value2 = 1

## 4. Resize the table / autofilter so the new row is included.
$lo.Resize($ws.Range("B9:G26"))

## 5. Fix up the conditional-formatting ranges that the engine does not
##    auto-shift on row insert.
$ws.Range("E9:E25").FormatConditions.Delete()
EOF_PLACEHOLDER
